$wb = $excel.ActiveWorkbook

# --- Sheet: Trends Status ---
$ws = $wb.Worksheets.Item("Trends Status")
$ws.Range("B2").Value = 99
$ws.Range("C2").Value = 65
$ws.Range("B3").Value = 105
$ws.Range("C3").Value = 81
$ws.Range("B4").Value = 420
$ws.Range("C4").Value = 284
$ws.Range("B5").Value = 186
$ws.Range("C5").Value = 300
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 190
$ws.Range("C7").Value = 14
$ws.Range("C8").Value = 12

# --- Sheet: Range Status ---
$ws = $wb.Worksheets.Item("Range Status")
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 298

# --- Sheet: Priority Status ---
$ws = $wb.Worksheets.Item("Priority Status")
$ws.Range("B2").Value = 159
$ws.Range("B3").Value = 385
$ws.Range("B4").Value = 402

# --- Sheet: Species qualification ---
$ws = $wb.Worksheets.Item("Species qualification")
$ws.Range("B2").Value = 946
$ws.Range("B3").Value = 533
$ws.Range("B4").Value = 673
$ws.Range("B5").Value = 946

# --- Sheet: SoIB-IUCN cross-tab ---
$ws = $wb.Worksheets.Item("SoIB-IUCN cross-tab")
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 46
$ws.Range("B6").Value = 75
$ws.Range("C6").Value = 327
$ws.Range("D6").Value = 385
$ws.Range("E6").Value = 787
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 5
$ws.Range("B8").Value = 159
$ws.Range("C8").Value = 402
$ws.Range("D8").Value = 385
$ws.Range("E8").Value = 946
